$d = $word.ActiveDocument

# 1) "Implementar el método solve() ..." paragraph: replace the closing
#    "...implementaríamos en una interface." with
#    "...implementaríamos en la clase concreta SilkRoadContest."
$d.Content.Find.Execute(
    "lo implementaríamos en una interface.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "lo implementaríamos en la clase concreta SilkRoadContest.", 2) | Out-Null

# 2) "Implementar el método simulate() ..." paragraph: replace the closing
#    "...esta lógica también se realiza en la interface del mini-ciclo 1."
#    with "...esta lógica también se realiza en la implementaríamos en la
#    clase concreta SilkRoadContest."
$d.Content.Find.Execute(
    "esta lógica también se realiza en la interface del mini-ciclo 1.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "esta lógica también se realiza en la implementaríamos en la clase concreta SilkRoadContest.", 2) | Out-Null

# 3) "Luiza -> 23 horas" -> "Luiza -> 26 horas"
$d.Content.Find.Execute(
    "Luiza -> 23 horas",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Luiza -> 26 horas", 2) | Out-Null

# 4) "Camilo -> 23 horas" -> "Camilo -> 26 horas"
$d.Content.Find.Execute(
    "Camilo -> 23 horas",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Camilo -> 26 horas", 2) | Out-Null

Write-Output "done"
